$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values updated per the "Updated cryptos list" GitHub Actions commit.
# Numeric-looking Price strings are quote-prefixed (leading "'") so Excel
# keeps storing them as text, matching the source inlineStr cells instead of
# silently coercing them into Number cells.
$ws.Range("D2").Value = "30.779.80"
$ws.Range("E2").Value = "  +2.62%  "
$ws.Range("D3").Value = "1.693.33"
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'222.27"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'30.89"
$ws.Range("E8").Value = "  +4.43%  "
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("D10").Value = "'0.0628"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "1.937.92"
$ws.Range("E12").Value = "  +3.55%  "
$ws.Range("D13").Value = "'10.65"
$ws.Range("E13").Value = "  +11.06%  "
$ws.Range("D14").Value = "'0.624"
$ws.Range("E14").Value = "  +8.11%  "
$ws.Range("D15").Value = "1.697.08"
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").Value = "30.775.20"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").Value = "'66.54"
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("D19").Value = "'249.91"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "'10.28"
$ws.Range("E22").Value = "  +4.48%  "
$ws.Range("D23").Value = "'4.31"
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").Value = "'157.88"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").Value = "'15.97"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "'0.0502"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("D31").Value = "'1.14"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("D33").Value = "1.520.69"
$ws.Range("E33").Value = "  +5.75%  "
$ws.Range("D34").Value = "'3.32"
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("E35").Value = "  +4.71%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("B37").Value = "Aave"
$ws.Range("C37").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D37").Value = "'81.49"
$ws.Range("E37").Value = "  +6.55%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0180"
$ws.Range("E38").Value = "  +4.96%  "
$ws.Range("E39").Value = "  +4.80%  "
$ws.Range("E40").Value = "  -5.75%  "
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").Value = "'0.856"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'0.998"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'52.58"
$ws.Range("E47").Value = "  -4.70%  "
$ws.Range("D48").Value = "1.830.16"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("D49").Value = "'5.45"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").Value = "'95.88"
$ws.Range("E50").Value = "  +6.05%  "
$ws.Range("D51").Value = "0.0₆0113"
$ws.Range("E51").Value = "  +2.51%  "
